# Switching to Summer time
# Shift all timestamps in column A (rows 2-97) forward by 3 days (DST change),
# and update the notified wind production values in column B (rows 2-93).
# Rows 94-97 in column B remain 0 (unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New literal timestamps for A2:A97 (explicit values avoid any rounding drift
# that could come from reading back a cell value and doing arithmetic on it)
$aValues = @(
    45744.01041666666,45744.02083333334,45744.03125,45744.04166666666,45744.05208333334,45744.0625,45744.07291666666,45744.08333333334,
    45744.09375,45744.10416666666,45744.11458333334,45744.125,45744.13541666666,45744.14583333334,45744.15625,45744.16666666666,
    45744.17708333334,45744.1875,45744.19791666666,45744.20833333334,45744.21875,45744.22916666666,45744.23958333334,45744.25,
    45744.26041666666,45744.27083333334,45744.28125,45744.29166666666,45744.30208333334,45744.3125,45744.32291666666,45744.33333333334,
    45744.34375,45744.35416666666,45744.36458333334,45744.375,45744.38541666666,45744.39583333334,45744.40625,45744.41666666666,
    45744.42708333334,45744.4375,45744.44791666666,45744.45833333334,45744.46875,45744.47916666666,45744.48958333334,45744.5,
    45744.51041666666,45744.52083333334,45744.53125,45744.54166666666,45744.55208333334,45744.5625,45744.57291666666,45744.58333333334,
    45744.59375,45744.60416666666,45744.61458333334,45744.625,45744.63541666666,45744.64583333334,45744.65625,45744.66666666666,
    45744.67708333334,45744.6875,45744.69791666666,45744.70833333334,45744.71875,45744.72916666666,45744.73958333334,45744.75,
    45744.76041666666,45744.77083333334,45744.78125,45744.79166666666,45744.80208333334,45744.8125,45744.82291666666,45744.83333333334,
    45744.84375,45744.85416666666,45744.86458333334,45744.875,45744.88541666666,45744.89583333334,45744.90625,45744.91666666666,
    45744.92708333334,45744.9375,45744.94791666666,45744.95833333334,45744.96875,45744.97916666666,45744.98958333334,45745
)

# New values for B2:B93
$bValues = @(
    1617,1641,1643,1651,1806,1814,1804,1803,2002,2014,2024,2038,2131,2141,2147,2155,
    2267,2268,2269,2261,2293,2294,2293,2291,2272,2265,2265,2266,2299,2299,2301,2302,
    2264,2267,2269,2273,2264,2266,2268,2269,2306,2308,2309,2311,2374,2375,2376,2377,
    2412,2413,2414,2414,2435,2434,2434,2434,2451,2450,2449,2448,2420,2417,2415,2412,
    2360,2358,2356,2353,2310,2309,2309,2308,2263,2259,2255,2251,2302,2299,2296,2293,
    2125,2122,2119,2116,2033,2031,2030,2028,2019,2020,2021,2022
)

$row = 2
foreach ($val in $aValues) {
    $ws.Cells.Item($row, 1).Value2 = $val
    $row++
}

$row = 2
foreach ($val in $bValues) {
    $ws.Cells.Item($row, 2).Value2 = $val
    $row++
}
